$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): extend sequence with P1=14, Q1=15 ---
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Copy formatting (bold font, border, centered) from O1 onto the new header cells
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Data rows 2..25 ---
# Swap the I/K and M/O values (1<->2) and append the two new columns P and Q (value 2)
for ($r = 2; $r -le 25; $r++) {
    $iVal = $ws.Cells.Item($r, 9).Value2
    $kVal = $ws.Cells.Item($r, 11).Value2
    $mVal = $ws.Cells.Item($r, 13).Value2
    $oVal = $ws.Cells.Item($r, 15).Value2

    $ws.Cells.Item($r, 9).Value = $kVal
    $ws.Cells.Item($r, 11).Value = $iVal
    $ws.Cells.Item($r, 13).Value = $oVal
    $ws.Cells.Item($r, 15).Value = $mVal

    $ws.Cells.Item($r, 16).Value = 2
    $ws.Cells.Item($r, 17).Value = 2
}
